$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AC2").Value = 19
$ws.Range("AD2").Value = 8.5
$ws.Range("AI2").Value = 9.5
$ws.Range("AK2").Value = 13
$ws.Range("AN2").Value = 7
$ws.Range("AO2").Value = 23
$ws.Range("AP2").Value = 26
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 3.75
$ws.Range("AY2").Value = 8
$ws.Range("BA2").Value = 21
$ws.Range("G2").Value = 5
$ws.Range("H2").Value = 4.33
$ws.Range("I2").Value = 1.6
$ws.Range("J2").Value = 4.75
$ws.Range("L2").Value = 2.1
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 19
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("AB3").Value = 29
$ws.Range("AE3").Value = 23
$ws.Range("AH3").Value = 26
$ws.Range("AJ3").Value = 29
$ws.Range("AK3").Value = 126
$ws.Range("AM3").Value = 51
$ws.Range("AQ3").Value = 13
$ws.Range("AT3").Value = 3.75
$ws.Range("AV3").Value = 67
$ws.Range("AX3").Value = 11
$ws.Range("AY3").Value = 51
$ws.Range("BA3").Value = 251
$ws.Range("G3").Value = 1.22
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 12
$ws.Range("L3").Value = 9.5
$ws.Range("M3").Value = 1.02
$ws.Range("N3").Value = 19
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 5.5
$ws.Range("Q3").Value = 1.5
$ws.Range("R3").Value = 2.5
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 3.75
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("W3").Value = 8.5
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 7.5
$ws.Range("AB5").Value = 41
$ws.Range("AC5").Value = 7
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 101
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 26
$ws.Range("AJ5").Value = 19
$ws.Range("AK5").Value = 67
$ws.Range("AL5").Value = 51
$ws.Range("AN5").Value = 3.4
$ws.Range("AO5").Value = 9.5
$ws.Range("AP5").Value = 26
$ws.Range("AR5").Value = 67
$ws.Range("AS5").Value = 251
$ws.Range("AT5").Value = 2.38
$ws.Range("AU5").Value = 10
$ws.Range("BA5").Value = 151
$ws.Range("BB5").Value = 201
$ws.Range("BC5").Value = 501
$ws.Range("BD5").Value = 126
$ws.Range("G5").Value = 1.7
$ws.Range("H5").Value = 3.5
$ws.Range("J5").Value = 2.4
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.63
$ws.Range("Q5").Value = 2.5
$ws.Range("R5").Value = 1.53
$ws.Range("S5").Value = 1.53
$ws.Range("T5").Value = 2.38
$ws.Range("U5").Value = 2.38
$ws.Range("V5").Value = 1.53
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 6.5
$ws.Range("AA6").Value = 17
$ws.Range("AL6").Value = 26
$ws.Range("BA6").Value = 51
$ws.Range("I6").Value = 3.5
$ws.Range("W6").Value = 8.5
$ws.Range("AA7").Value = 34
$ws.Range("AD7").Value = 7
$ws.Range("AH7").Value = 8
$ws.Range("AI7").Value = 9
$ws.Range("AM7").Value = 23
$ws.Range("AO7").Value = 23
$ws.Range("AP7").Value = 29
$ws.Range("AY7").Value = 9
$ws.Range("G7").Value = 4.2
$ws.Range("I7").Value = 1.75
$ws.Range("J7").Value = 4.5
$ws.Range("L7").Value = 2.38
$ws.Range("X7").Value = 23
$ws.Range("Y7").Value = 15
$ws.Range("AB8").Value = 34
$ws.Range("AC8").Value = 11
$ws.Range("AG8").Value = 201
$ws.Range("AH8").Value = 8.5
$ws.Range("AI8").Value = 11
$ws.Range("AJ8").Value = 9
$ws.Range("AM8").Value = 26
$ws.Range("AP8").Value = 26
$ws.Range("AR8").Value = 81
$ws.Range("AT8").Value = 2.75
$ws.Range("AX8").Value = 4.33
$ws.Range("AY8").Value = 12
$ws.Range("BC8").Value = 151
$ws.Range("G8").Value = 3.2
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 2.2
$ws.Range("J8").Value = 3.75
$ws.Range("L8").Value = 2.88
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.93
$ws.Range("R8").Value = 1.93
$ws.Range("S8").Value = 1.4
$ws.Range("T8").Value = 2.75
$ws.Range("U8").Value = 1.73
$ws.Range("V8").Value = 2
$ws.Range("W8").Value = 10
$ws.Range("X8").Value = 17
$ws.Range("Z8").Value = 34
$ws.Range("AB9").Value = 34
$ws.Range("AC9").Value = 9
$ws.Range("AE9").Value = 17
$ws.Range("AG9").Value = 351
$ws.Range("AH9").Value = 8.5
$ws.Range("AJ9").Value = 12
$ws.Range("AL9").Value = 29
$ws.Range("AM9").Value = 41
$ws.Range("AO9").Value = 13
$ws.Range("AS9").Value = 201
$ws.Range("AT9").Value = 2.63
$ws.Range("AU9").Value = 8.5
$ws.Range("AZ9").Value = 29
$ws.Range("BA9").Value = 67
$ws.Range("BC9").Value = 251
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 3.2
$ws.Range("J9").Value = 3
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 9
$ws.Range("O9").Value = 1.4
$ws.Range("P9").Value = 3
$ws.Range("Q9").Value = 2.2
$ws.Range("R9").Value = 1.67
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 2.63
$ws.Range("U9").Value = 1.83
$ws.Range("V9").Value = 1.83
$ws.Range("W9").Value = 7
$ws.Range("X9").Value = 10
$ws.Range("Y9").Value = 9.5
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 13
$ws.Range("BC10").Value = 251
$ws.Range("G10").Value = 1.95
$ws.Range("I10").Value = 4
$ws.Range("J10").Value = 2.63
$ws.Range("L10").Value = 4.5
$ws.Range("Z10").Value = 17
$ws.Range("AC11").Value = 12
$ws.Range("AE11").Value = 15
$ws.Range("AF11").Value = 51
$ws.Range("AH11").Value = 12
$ws.Range("AL11").Value = 34
$ws.Range("AO11").Value = 10
$ws.Range("AQ11").Value = 34
$ws.Range("AU11").Value = 8
$ws.Range("AZ11").Value = 29
$ws.Range("BC11").Value = 201
$ws.Range("H11").Value = 3.6
$ws.Range("J11").Value = 2.5
$ws.Range("K11").Value = 2.2
$ws.Range("L11").Value = 4.5
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("O11").Value = 1.29
$ws.Range("P11").Value = 3.75
$ws.Range("Q11").Value = 1.93
$ws.Range("R11").Value = 1.93
$ws.Range("U11").Value = 1.73
$ws.Range("V11").Value = 2
$ws.Range("X11").Value = 9
$ws.Range("AA12").Value = 13.5
$ws.Range("AB12").Value = 28
$ws.Range("AC12").Value = 9.5
$ws.Range("AD12").Value = 7
$ws.Range("AE12").Value = 17
$ws.Range("AF12").Value = 90
$ws.Range("AG12").Value = 700
$ws.Range("AH12").Value = 12.5
$ws.Range("AI12").Value = 29
$ws.Range("AJ12").Value = 16.5
$ws.Range("AK12").Value = 100
$ws.Range("AL12").Value = 55
$ws.Range("AM12").Value = 55
$ws.Range("AO12").Value = 8
$ws.Range("AP12").Value = 18
$ws.Range("AQ12").Value = 27
$ws.Range("AT12").Value = 2.47
$ws.Range("AU12").Value = 7.7
$ws.Range("AV12").Value = 80
$ws.Range("AY12").Value = 30
$ws.Range("BA12").Value = 200
$ws.Range("BB12").Value = 250
$ws.Range("BC12").Value = 500
$ws.Range("G12").Value = 1.65
$ws.Range("H12").Value = 3.6
$ws.Range("I12").Value = 4.9
$ws.Range("K12").Value = 2.1
$ws.Range("L12").Value = 5.2
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("O12").Value = 1.29
$ws.Range("P12").Value = 2.95
$ws.Range("Q12").Value = 1.87
$ws.Range("R12").Value = 1.75
$ws.Range("S12").Value = 1.4
$ws.Range("T12").Value = 2.52
$ws.Range("U12").Value = 1.85
$ws.Range("V12").Value = 1.75
$ws.Range("W12").Value = 6.4
$ws.Range("X12").Value = 7.4
$ws.Range("N13").Value = 8.85
$ws.Range("AC14").Value = 8
$ws.Range("AD14").Value = 6.6
$ws.Range("AH14").Value = 10.5
$ws.Range("AI14").Value = 16
$ws.Range("AO14").Value = 12
$ws.Range("AQ14").Value = 45
$ws.Range("AY14").Value = 15
$ws.Range("AZ14").Value = 20
$ws.Range("H14").Value = 3.25
$ws.Range("I14").Value = 2.75
$ws.Range("J14").Value = 2.87
$ws.Range("N14").Value = 8
$ws.Range("W14").Value = 10.25
